$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4000963.8
$ws.Range("I33").Value = 5882574
$ws.Range("K33").Value = 5882574
$ws.Range("M33").Value = -5882345
$ws.Range("H112").Value = 78584.38
$ws.Range("J112").Value = 78584.38
$ws.Range("L112").Value = 235753.14
$ws.Range("N112").Value = -237969.14
$ws.Range("H131").Value = 32800
$ws.Range("I131").Value = 2550
$ws.Range("J131").Value = 41442.855
$ws.Range("K131").Value = 7650
$ws.Range("L131").Value = 124328.565
$ws.Range("M131").Value = -2610
$ws.Range("N131").Value = -134408.565

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3292.6
$ws.Range("I61").Value = 2914.4
$ws.Range("K61").Value = 2914.4
$ws.Range("M61").Value = -2702.4
$ws.Range("H110").Value = 11407.208
$ws.Range("I110").Value = 16735.334
$ws.Range("J110").Value = 6079.0835
$ws.Range("K110").Value = 16735.334
$ws.Range("L110").Value = 6079.0835
$ws.Range("M110").Value = -14690.334
$ws.Range("N110").Value = -10169.0835
$ws.Range("H113").Value = 63157
$ws.Range("J113").Value = 63157
$ws.Range("L113").Value = 63157
$ws.Range("N113").Value = -71835
$ws.Range("H132").Value = 5769.778
$ws.Range("I132").Value = 5679.4
$ws.Range("J132").Value = 5882.75
$ws.Range("K132").Value = 17038.2
$ws.Range("L132").Value = 17648.25
$ws.Range("M132").Value = -14508.2
$ws.Range("N132").Value = -22708.25
$ws.Range("H136").Value = 3292.6
$ws.Range("I136").Value = 2914.4
$ws.Range("K136").Value = 8743.200000000001
$ws.Range("M136").Value = -6193.200000000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5200
$ws.Range("I4").Value = 3600
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 3600
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = -3488
$ws.Range("N4").Value = -10224
$ws.Range("H31").Value = 57915.473
$ws.Range("I31").Value = 102148.3
$ws.Range("J31").Value = 8767.888999999999
$ws.Range("K31").Value = 102148.3
$ws.Range("L31").Value = 8767.888999999999
$ws.Range("M31").Value = -101853.3
$ws.Range("N31").Value = -9357.888999999999
$ws.Range("H34").Value = 57915.473
$ws.Range("I34").Value = 102148.3
$ws.Range("J34").Value = 8767.888999999999
$ws.Range("K34").Value = 102148.3
$ws.Range("L34").Value = 8767.888999999999
$ws.Range("M34").Value = -101946.3
$ws.Range("N34").Value = -9171.888999999999
$ws.Range("H58").Value = 2061
$ws.Range("I58").Value = 2118.1875
$ws.Range("J58").Value = 1832.25
$ws.Range("K58").Value = 2118.1875
$ws.Range("L58").Value = 1832.25
$ws.Range("M58").Value = -1915.1875
$ws.Range("N58").Value = -2238.25
$ws.Range("H132").Value = 2549.9312
$ws.Range("I132").Value = 2193.6155
$ws.Range("K132").Value = 6580.8465
$ws.Range("M132").Value = -4050.8465
$ws.Range("H134").Value = 26272.143
$ws.Range("I134").Value = 21281
$ws.Range("J134").Value = 38750
$ws.Range("K134").Value = 63843
$ws.Range("L134").Value = 116250
$ws.Range("M134").Value = -61308
$ws.Range("N134").Value = -121320
$ws.Range("H136").Value = 2061
$ws.Range("I136").Value = 2118.1875
$ws.Range("J136").Value = 1832.25
$ws.Range("K136").Value = 6354.5625
$ws.Range("L136").Value = 5496.75
$ws.Range("M136").Value = -3804.5625
$ws.Range("N136").Value = -10596.75

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 150.66667
$ws.Range("I17").Value = 139.33333
$ws.Range("J17").Value = 184.66667
$ws.Range("K17").Value = 417.99999
$ws.Range("L17").Value = 554.00001
$ws.Range("M17").Value = -248.99999
$ws.Range("N17").Value = -892.00001
$ws.Range("H40").Value = 312.90475
$ws.Range("I40").Value = 198.45454
$ws.Range("K40").Value = 793.81816
$ws.Range("M40").Value = -724.81816
$ws.Range("H107").Value = 1093.6
$ws.Range("I107").Value = 270
$ws.Range("K107").Value = 810
$ws.Range("M107").Value = 1110
$ws.Range("H122").Value = 3314.6
$ws.Range("I122").Value = 3036
$ws.Range("J122").Value = 3500.3333
$ws.Range("K122").Value = 27324
$ws.Range("L122").Value = 31502.9997
$ws.Range("M122").Value = -24874
$ws.Range("N122").Value = -36402.9997
$ws.Range("H132").Value = 1299.9333
$ws.Range("I132").Value = 1045.3636
$ws.Range("K132").Value = 9408.2724
$ws.Range("M132").Value = -6878.2724

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = $null
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = $null
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = $null
$ws.Range("H113").Value = 1779.8
$ws.Range("I113").Value = 1779.8
$ws.Range("K113").Value = 1779.8
$ws.Range("M113").Value = 390.2

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1158.3334
$ws.Range("I22").Value = 562.5
$ws.Range("K22").Value = 562.5
$ws.Range("M22").Value = -267.5
$ws.Range("H27").Value = 1158.3334
$ws.Range("I27").Value = 562.5
$ws.Range("K27").Value = 562.5
$ws.Range("M27").Value = -455.5
$ws.Range("H55").Value = 287.0909
$ws.Range("I55").Value = 293.42856
$ws.Range("J55").Value = 276
$ws.Range("K55").Value = 293.42856
$ws.Range("L55").Value = 276
$ws.Range("M55").Value = -120.42856
$ws.Range("N55").Value = -622
$ws.Range("H61").Value = 3986.0715
$ws.Range("I61").Value = 4216.6665
$ws.Range("J61").Value = 2602.5
$ws.Range("K61").Value = 4216.6665
$ws.Range("L61").Value = 2602.5
$ws.Range("M61").Value = -4014.6665
$ws.Range("N61").Value = -3006.5
$ws.Range("H113").Value = 3986.0715
$ws.Range("I113").Value = 4216.6665
$ws.Range("J113").Value = 2602.5
$ws.Range("K113").Value = 4216.6665
$ws.Range("L113").Value = 2602.5
$ws.Range("M113").Value = -2046.6665
$ws.Range("N113").Value = -6942.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 7500
$ws.Range("I7").Value = 7500
$ws.Range("K7").Value = 7500
$ws.Range("M7").Value = -7387
$ws.Range("H9").Value = 49998.5
$ws.Range("I9").Value = 49998.5
$ws.Range("K9").Value = 49998.5
$ws.Range("M9").Value = -49858.5
$ws.Range("H28").Value = 14180
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 14180
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 14180
$ws.Range("M28").Value = $null
$ws.Range("N28").Value = -14876
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = $null
$ws.Range("H62").Value = 22388
$ws.Range("I62").Value = 30649.285
$ws.Range("J62").Value = 12749.833
$ws.Range("K62").Value = 30649.285
$ws.Range("L62").Value = 12749.833
$ws.Range("M62").Value = -30025.285
$ws.Range("N62").Value = -13997.833
$ws.Range("H65").Value = 22388
$ws.Range("I65").Value = 30649.285
$ws.Range("J65").Value = 12749.833
$ws.Range("K65").Value = 153246.425
$ws.Range("L65").Value = 63749.165
$ws.Range("M65").Value = -150126.425
$ws.Range("N65").Value = -69989.16500000001
$ws.Range("H122").Value = 2803.4443
$ws.Range("I122").Value = 2792.76
$ws.Range("J122").Value = 2827.7273
$ws.Range("K122").Value = 8378.280000000001
$ws.Range("L122").Value = 8483.1819
$ws.Range("M122").Value = -5928.280000000001
$ws.Range("N122").Value = -13383.1819
$ws.Range("H126").Value = 4833.3335
$ws.Range("I126").Value = 4250
$ws.Range("K126").Value = 12750
$ws.Range("M126").Value = -10280
$ws.Range("H132").Value = 3618.3547
$ws.Range("I132").Value = 3547.1785
$ws.Range("J132").Value = 4282.6665
$ws.Range("K132").Value = 4282.6665
$ws.Range("L132").Value = 12847.9995
$ws.Range("M132").Value = -8111.5355
$ws.Range("N132").Value = -17907.9995
